# Fruta / hortaliza, semanal
# Insert 3 new weekly observation rows before the current row 535, pushing the
# existing rows 535:541 down to 538:544, then populate the 3 new rows with
# their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 535 (shifts old 535:541 -> 538:544, carrying
# formatting such as the column D date style down with them).
$ws.Range("A535:R537").EntireRow.Insert()

# --- Row 535 ---
$ws.Cells.Item(535, 1).Value = 6
$ws.Cells.Item(535, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(535, 3).Value = "Metropolitana"
$ws.Cells.Item(535, 4).Value = 44595
$ws.Cells.Item(535, 5).Value = 13
$ws.Cells.Item(535, 6).Value = 100112017
$ws.Cells.Item(535, 7).Value = "Apio"
$ws.Cells.Item(535, 8).Value = "Americana (o)"
$ws.Cells.Item(535, 9).Value = "Primera"
$ws.Cells.Item(535, 10).Value = 1500
$ws.Cells.Item(535, 11).Value = 6000
$ws.Cells.Item(535, 12).Value = 7000
$ws.Cells.Item(535, 13).Value = 6400
$ws.Cells.Item(535, 14).Value = "`$/docena de matas"
$ws.Cells.Item(535, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(535, 16).Value = 1067
$ws.Cells.Item(535, 17).Value = 6
$ws.Cells.Item(535, 18).Value = "Hortaliza"

# --- Row 536 ---
$ws.Cells.Item(536, 1).Value = 6
$ws.Cells.Item(536, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(536, 3).Value = "Metropolitana"
$ws.Cells.Item(536, 4).Value = 44595
$ws.Cells.Item(536, 5).Value = 13
$ws.Cells.Item(536, 6).Value = 100112017
$ws.Cells.Item(536, 7).Value = "Apio"
$ws.Cells.Item(536, 8).Value = "Americana (o)"
$ws.Cells.Item(536, 9).Value = "Segunda"
$ws.Cells.Item(536, 10).Value = 230
$ws.Cells.Item(536, 11).Value = 4000
$ws.Cells.Item(536, 12).Value = 4000
$ws.Cells.Item(536, 13).Value = 4000
$ws.Cells.Item(536, 14).Value = "`$/docena de matas"
$ws.Cells.Item(536, 15).Value = "Región Metropolitana"
$ws.Cells.Item(536, 16).Value = 667
$ws.Cells.Item(536, 17).Value = 6
$ws.Cells.Item(536, 18).Value = "Hortaliza"

# --- Row 537 ---
$ws.Cells.Item(537, 1).Value = 6
$ws.Cells.Item(537, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(537, 3).Value = "Metropolitana"
$ws.Cells.Item(537, 4).Value = 44595
$ws.Cells.Item(537, 5).Value = 13
$ws.Cells.Item(537, 6).Value = 100112017
$ws.Cells.Item(537, 7).Value = "Apio"
$ws.Cells.Item(537, 8).Value = "Americana (o)"
$ws.Cells.Item(537, 9).Value = "Segunda"
$ws.Cells.Item(537, 10).Value = 800
$ws.Cells.Item(537, 11).Value = 4500
$ws.Cells.Item(537, 12).Value = 5000
$ws.Cells.Item(537, 13).Value = 4812
$ws.Cells.Item(537, 14).Value = "`$/docena de matas"
$ws.Cells.Item(537, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(537, 16).Value = 802
$ws.Cells.Item(537, 17).Value = 6
$ws.Cells.Item(537, 18).Value = "Hortaliza"
